# Update "想去人数" (wanting-to-go count) figures for two events that
# appear on both the "展览" sheet and the consolidated "全部类型" sheet.
#
#   展览   row 3 (F3): 1239 -> 1241
#   展览   row 7 (F7): 161  -> 162
#   全部类型 row 3 (F3): 1239 -> 1241
#   全部类型 row 7 (F7): 161  -> 162

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1241
    $ws.Range("F7").Value = 162
}
